$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: swap date with row 4, update volume
$ws.Range("D3").Value = 44250
$ws.Range("M3").Value = 200

# Row 4: swap date with row 3, update volume
$ws.Range("D4").Value = 44253
$ws.Range("M4").Value = 160

# Row 5: takes on row 6's former values
$ws.Range("D5").Value = 45072
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("R5").Value = "Provincia de Chacabuco"
$ws.Range("S5").Value = 889

# Row 6: takes on row 7's former values (D6/L6 unchanged)
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 17000
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 944

# Row 7: takes on row 5's former (original) values
$ws.Range("D7").Value = 44257
$ws.Range("L7").Value = "Primera"
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 806
